# Update "want to go" counts (column F) on each sheet to the freshly
# scraped values (gh-pages output regenerated at commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 2680
$ws.Range("F3").Value = 1032
$ws.Range("F4").Value = 19236
$ws.Range("F6").Value = 2156
$ws.Range("F7").Value = 727
$ws.Range("F8").Value = 609
$ws.Range("F9").Value = 406
$ws.Range("F10").Value = 666
$ws.Range("F11").Value = 223
$ws.Range("F12").Value = 237
$ws.Range("F14").Value = 348
$ws.Range("F16").Value = 245
$ws.Range("F20").Value = 17
$ws.Range("F21").Value = 87

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 183
$ws.Range("F3").Value = 33
$ws.Range("F7").Value = 264
$ws.Range("F13").Value = 82
$ws.Range("F15").Value = 56
$ws.Range("F18").Value = 4

# Sheet 3: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5963
$ws.Range("F3").Value = 621
$ws.Range("F4").Value = 577

# Sheet 4: 全部类型 (All Types) - union of the above three sheets
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 5963
$ws.Range("F3").Value = 621
$ws.Range("F4").Value = 577
$ws.Range("F5").Value = 183
$ws.Range("F6").Value = 33
$ws.Range("F7").Value = 2680
$ws.Range("F8").Value = 1032
$ws.Range("F9").Value = 19237
$ws.Range("F14").Value = 264
$ws.Range("F15").Value = 2156
$ws.Range("F16").Value = 727
$ws.Range("F18").Value = 609
$ws.Range("F19").Value = 406
$ws.Range("F20").Value = 666
$ws.Range("F21").Value = 223
$ws.Range("F22").Value = 237
$ws.Range("F27").Value = 348
$ws.Range("F30").Value = 245
$ws.Range("F31").Value = 82
$ws.Range("F35").Value = 56
$ws.Range("F39").Value = 4
$ws.Range("F40").Value = 17
$ws.Range("F45").Value = 87
